$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 99
$ws.Range("I2").Value = 302
$ws.Range("J2").Value = 1144
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 324
$ws.Range("M2").Value = 19
$ws.Range("N2").Value = 240
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 128
$ws.Range("T2").Value = 202
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 1830
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1904
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 8
